# "4.0.3 model and data"
#
# The "Boolean" sheet's InputData pathname list splits two aggregate
# transportation CSVs (trans/BVTQaZ/BVTQaZ.csv and trans/VTQaZ/VTQaZ.csv)
# into six mode-specific CSVs each (LDVs, HDVs, aircraft, rail, ships,
# motorbikes). The "About" sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boolean")

# --- Split "trans/BVTQaZ/BVTQaZ.csv" (row 17) into six rows ---------------
# Insert 5 new blank rows above row 17 so rows 17-22 are available, then
# write the six replacement values (overwriting the old BVTQaZ.csv row that
# shifted down to row 22 along with the inserted blanks).
$ws.Range("A17:A21").EntireRow.Insert()

$ws.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$ws.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$ws.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$ws.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$ws.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$ws.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the above, the unaffected rows are now:
#   23 trans/BVTStL/BVTStL.csv
#   24 trans/PVTStL/PVTStL.csv
#   25 trans/SRPbVT/SRPbVT.csv
#   26 trans/VTQaZ/VTQaZ.csv            <- split into six rows next
#   27 trans/VTStFES/VTStFES.csv

# --- Split "trans/VTQaZ/VTQaZ.csv" (now row 26) into six rows -------------
$ws.Range("A26:A30").EntireRow.Insert()

$ws.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$ws.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$ws.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$ws.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$ws.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$ws.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Row 32 now holds trans/VTStFES/VTStFES.csv (unchanged), followed by six
# trailing blank rows (33-38) left over from editing/scrolling in the app.
$ws.Range("A33:A38").EntireRow.Insert()

# --- View state tweaks ------------------------------------------------------
$ws.Application.GoTo($ws.Range("A32"), $true)

$wsInteger = $wb.Worksheets.Item("Integer")
$wsInteger.Range("A13").Select() | Out-Null

# "About" becomes the selected/active tab (was "Integer").
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
